# Recomputed TPM-based NATMI ligand/receptor metrics for Hbegf-Cd44 (YoungD7).
# Updates columns G,H,I,J (ligand avg/total expr + derived specificity) and
# M,N,O,P (receptor avg/total expr + derived specificity) plus the dependent
# edge-weight columns Q,R,S,T for rows 2-17, per the new TPM recalculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 8.082188
$ws.Range("H2").Value = 24.246564
$ws.Range("I2").Value = 0.2755301789948819
$ws.Range("J2").Value = 0.2755301789948819
$ws.Range("M2").Value = 24.576554
$ws.Range("N2").Value = 73.729662
$ws.Range("O2").Value = 0.07553767049546639
$ws.Range("P2").Value = 0.07553767049546638
$ws.Range("Q2").Value = 198.632329820152
$ws.Range("R2").Value = 1787.690968381368
$ws.Range("S2").Value = 0.02081290787247226
$ws.Range("T2").Value = 0.02081290787247226

# Row 3
$ws.Range("G3").Value = 8.082188
$ws.Range("H3").Value = 24.246564
$ws.Range("I3").Value = 0.2755301789948819
$ws.Range("J3").Value = 0.2755301789948819
$ws.Range("O3").Value = 0.359764849016532
$ws.Range("P3").Value = 0.359764849016532
$ws.Range("Q3").Value = 946.0303670846974
$ws.Range("R3").Value = 8514.273303762278
$ws.Range("S3").Value = 0.09912607324559171
$ws.Range("T3").Value = 0.09912607324559172

# Row 4
$ws.Range("G4").Value = 8.082188
$ws.Range("H4").Value = 24.246564
$ws.Range("I4").Value = 0.2755301789948819
$ws.Range("J4").Value = 0.2755301789948819
$ws.Range("M4").Value = 55.68784966666667
$ws.Range("N4").Value = 167.063549
$ws.Range("O4").Value = 0.1711603033819035
$ws.Range("P4").Value = 0.1711603033819035
$ws.Range("Q4").Value = 450.0796703217374
$ws.Range("R4").Value = 4050.717032895637
$ws.Range("S4").Value = 0.04715982902763416
$ws.Range("T4").Value = 0.04715982902763417

# Row 5
$ws.Range("G5").Value = 8.082188
$ws.Range("H5").Value = 24.246564
$ws.Range("I5").Value = 0.2755301789948819
$ws.Range("J5").Value = 0.2755301789948819
$ws.Range("M5").Value = 128.0392633333333
$ws.Range("N5").Value = 384.11779
$ws.Range("O5").Value = 0.3935371771060981
$ws.Range("P5").Value = 0.3935371771060981
$ws.Range("Q5").Value = 1034.837397641507
$ws.Range("R5").Value = 9313.536578773561
$ws.Range("S5").Value = 0.1084313688491837
$ws.Range("T5").Value = 0.1084313688491838

# Row 6
$ws.Range("I6").Value = 0.4533445791334642
$ws.Range("J6").Value = 0.4533445791334642
$ws.Range("M6").Value = 24.576554
$ws.Range("N6").Value = 73.729662
$ws.Range("O6").Value = 0.07553767049546639
$ws.Range("P6").Value = 0.07553767049546638
$ws.Range("Q6").Value = 326.8204241477627
$ws.Range("R6").Value = 2941.383817329864
$ws.Range("S6").Value = 0.0342445934394895
$ws.Range("T6").Value = 0.0342445934394895

# Row 7
$ws.Range("I7").Value = 0.4533445791334642
$ws.Range("J7").Value = 0.4533445791334642
$ws.Range("O7").Value = 0.359764849016532
$ws.Range("P7").Value = 0.359764849016532
$ws.Range("S7").Value = 0.163097444064414
$ws.Range("T7").Value = 0.163097444064414

# Row 8
$ws.Range("I8").Value = 0.4533445791334642
$ws.Range("J8").Value = 0.4533445791334642
$ws.Range("M8").Value = 55.68784966666667
$ws.Range("N8").Value = 167.063549
$ws.Range("O8").Value = 0.1711603033819035
$ws.Range("P8").Value = 0.1711603033819035
$ws.Range("Q8").Value = 740.5402176373809
$ws.Range("R8").Value = 6664.861958736428
$ws.Range("S8").Value = 0.07759459570102509
$ws.Range("T8").Value = 0.07759459570102509

# Row 9
$ws.Range("I9").Value = 0.4533445791334642
$ws.Range("J9").Value = 0.4533445791334642
$ws.Range("M9").Value = 128.0392633333333
$ws.Range("N9").Value = 384.11779
$ws.Range("O9").Value = 0.3935371771060981
$ws.Range("P9").Value = 0.3935371771060981
$ws.Range("Q9").Value = 1702.673464724431
$ws.Range("R9").Value = 15324.06118251988
$ws.Range("S9").Value = 0.1784079459285356
$ws.Range("T9").Value = 0.1784079459285356

# Row 10
$ws.Range("G10").Value = 5.789497666666667
$ws.Range("H10").Value = 17.368493
$ws.Range("I10").Value = 0.1973699855023315
$ws.Range("J10").Value = 0.1973699855023315
$ws.Range("M10").Value = 24.576554
$ws.Range("N10").Value = 73.729662
$ws.Range("O10").Value = 0.07553767049546639
$ws.Range("P10").Value = 0.07553767049546638
$ws.Range("Q10").Value = 142.2859020377073
$ws.Range("R10").Value = 1280.573118339366
$ws.Range("S10").Value = 0.01490886893057009
$ws.Range("T10").Value = 0.01490886893057009

# Row 11
$ws.Range("G11").Value = 5.789497666666667
$ws.Range("H11").Value = 17.368493
$ws.Range("I11").Value = 0.1973699855023315
$ws.Range("J11").Value = 0.1973699855023315
$ws.Range("O11").Value = 0.359764849016532
$ws.Range("P11").Value = 0.359764849016532
$ws.Range("Q11").Value = 677.6680526155375
$ws.Range("R11").Value = 6099.012473539838
$ws.Range("S11").Value = 0.07100678303464139
$ws.Range("T11").Value = 0.0710067830346414

# Row 12
$ws.Range("G12").Value = 5.789497666666667
$ws.Range("H12").Value = 17.368493
$ws.Range("I12").Value = 0.1973699855023315
$ws.Range("J12").Value = 0.1973699855023315
$ws.Range("M12").Value = 55.68784966666667
$ws.Range("N12").Value = 167.063549
$ws.Range("O12").Value = 0.1711603033819035
$ws.Range("P12").Value = 0.1711603033819035
$ws.Range("Q12").Value = 322.4046757068508
$ws.Range("R12").Value = 2901.642081361657
$ws.Range("S12").Value = 0.03378190659706096
$ws.Range("T12").Value = 0.03378190659706096

# Row 13
$ws.Range("G13").Value = 5.789497666666667
$ws.Range("H13").Value = 17.368493
$ws.Range("I13").Value = 0.1973699855023315
$ws.Range("J13").Value = 0.1973699855023315
$ws.Range("M13").Value = 128.0392633333333
$ws.Range("N13").Value = 384.11779
$ws.Range("O13").Value = 0.3935371771060981
$ws.Range("P13").Value = 0.3935371771060981
$ws.Range("Q13").Value = 741.2830163100523
$ws.Range("R13").Value = 6671.547146790471
$ws.Range("S13").Value = 0.07767242694005903
$ws.Range("T13").Value = 0.07767242694005905

# Row 14
$ws.Range("G14").Value = 2.163479333333334
$ws.Range("H14").Value = 6.490438
$ws.Range("I14").Value = 0.0737552563693224
$ws.Range("J14").Value = 0.0737552563693224
$ws.Range("M14").Value = 24.576554
$ws.Range("N14").Value = 73.729662
$ws.Range("O14").Value = 0.07553767049546639
$ws.Range("P14").Value = 0.07553767049546638
$ws.Range("Q14").Value = 53.17086666355068
$ws.Range("R14").Value = 478.537799971956
$ws.Range("S14").Value = 0.005571300252934524
$ws.Range("T14").Value = 0.005571300252934523

# Row 15
$ws.Range("G15").Value = 2.163479333333334
$ws.Range("H15").Value = 6.490438
$ws.Range("I15").Value = 0.0737552563693224
$ws.Range("J15").Value = 0.0737552563693224
$ws.Range("O15").Value = 0.359764849016532
$ws.Range("P15").Value = 0.359764849016532
$ws.Range("Q15").Value = 253.2380028642602
$ws.Range("R15").Value = 2279.142025778342
$ws.Range("S15").Value = 0.02653454867188488
$ws.Range("T15").Value = 0.02653454867188488

# Row 16
$ws.Range("G16").Value = 2.163479333333334
$ws.Range("H16").Value = 6.490438
$ws.Range("I16").Value = 0.0737552563693224
$ws.Range("J16").Value = 0.0737552563693224
$ws.Range("M16").Value = 55.68784966666667
$ws.Range("N16").Value = 167.063549
$ws.Range("O16").Value = 0.1711603033819035
$ws.Range("P16").Value = 0.1711603033819035
$ws.Range("Q16").Value = 120.4795118716069
$ws.Range("R16").Value = 1084.315606844462
$ws.Range("S16").Value = 0.01262397205618329
$ws.Range("T16").Value = 0.01262397205618329

# Row 17
$ws.Range("G17").Value = 2.163479333333334
$ws.Range("H17").Value = 6.490438
$ws.Range("I17").Value = 0.0737552563693224
$ws.Range("J17").Value = 0.0737552563693224
$ws.Range("M17").Value = 128.0392633333333
$ws.Range("N17").Value = 384.11779
$ws.Range("O17").Value = 0.3935371771060981
$ws.Range("P17").Value = 0.3935371771060981
$ws.Range("Q17").Value = 277.0103000768912
$ws.Range("R17").Value = 2493.09270069202
$ws.Range("S17").Value = 0.0290254353883197
$ws.Range("T17").Value = 0.0290254353883197
